$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 145.64285
$ws.Range("I9").Value = 146.9
$ws.Range("J9").Value = 142.5
$ws.Range("K9").Value = 146.9
$ws.Range("L9").Value = 142.5
$ws.Range("M9").Value = 22.09999999999999
$ws.Range("N9").Value = -480.5

$ws.Range("H32").Value = 2529.6316
$ws.Range("I32").Value = 3355.7144
$ws.Range("J32").Value = 2047.75
$ws.Range("K32").Value = 3355.7144
$ws.Range("L32").Value = 2047.75
$ws.Range("M32").Value = -3029.7144
$ws.Range("N32").Value = -2699.75

$ws.Range("H39").Value = 583.7143
$ws.Range("I39").Value = 179
$ws.Range("J39").Value = 887.25
$ws.Range("K39").Value = 537
$ws.Range("L39").Value = 2661.75
$ws.Range("M39").Value = -241
$ws.Range("N39").Value = -3253.75

$ws.Range("H69").Value = 4440
$ws.Range("I69").Value = 2000
$ws.Range("J69").Value = 5050
$ws.Range("K69").Value = 6000
$ws.Range("L69").Value = 15150
$ws.Range("M69").Value = -5126
$ws.Range("N69").Value = -16898

$ws.Range("H70").Value = 1171.4286
$ws.Range("I70").Value = 700
$ws.Range("J70").Value = 1250
$ws.Range("K70").Value = 2100
$ws.Range("L70").Value = 3750
$ws.Range("M70").Value = -1830
$ws.Range("N70").Value = -4290

$ws.Range("H72").Value = 4440
$ws.Range("I72").Value = 2000
$ws.Range("J72").Value = 5050
$ws.Range("K72").Value = 18000
$ws.Range("L72").Value = 45450
$ws.Range("M72").Value = -13632
$ws.Range("N72").Value = -54186

$ws.Range("H73").Value = 1171.4286
$ws.Range("I73").Value = 700
$ws.Range("J73").Value = 1250
$ws.Range("K73").Value = 2100
$ws.Range("L73").Value = 3750
$ws.Range("M73").Value = -1164
$ws.Range("N73").Value = -5622

$ws.Range("H112").Value = 1365.4546
$ws.Range("J112").Value = 1496.8422
$ws.Range("L112").Value = 4490.5266
$ws.Range("N112").Value = -6706.5266

$ws.Range("H129").Value = 434.2857
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws.Range("H132").Value = 6527.5884
$ws.Range("I132").Value = 7083.5
$ws.Range("J132").Value = 3933.3333
$ws.Range("K132").Value = 21250.5
$ws.Range("L132").Value = 11799.9999
$ws.Range("M132").Value = -18720.5
$ws.Range("N132").Value = -16859.9999

$ws.Range("H138").Value = 2366.724
$ws.Range("I138").Value = 1946.381
$ws.Range("K138").Value = 5839.143
$ws.Range("M138").Value = -699.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5530.38
$ws.Range("I32").Value = 5449.367
$ws.Range("K32").Value = 5449.367
$ws.Range("M32").Value = -5162.367

$ws.Range("H61").Value = 4935.8696
$ws.Range("I61").Value = 5334.294
$ws.Range("K61").Value = 5334.294
$ws.Range("M61").Value = -5122.294

$ws.Range("H122").Value = 1482.7333
$ws.Range("I122").Value = 1392.6666
$ws.Range("J122").Value = 1843
$ws.Range("K122").Value = 4177.9998
$ws.Range("L122").Value = 5529
$ws.Range("M122").Value = -1727.9998
$ws.Range("N122").Value = -10429

$ws.Range("H132").Value = 60587.234
$ws.Range("I132").Value = 1289.5454
$ws.Range("J132").Value = 169299.67
$ws.Range("K132").Value = 3868.6362
$ws.Range("L132").Value = 507899.01
$ws.Range("M132").Value = -1338.6362
$ws.Range("N132").Value = -512959.01

$ws.Range("H136").Value = 4935.8696
$ws.Range("I136").Value = 5334.294
$ws.Range("K136").Value = 16002.882
$ws.Range("M136").Value = -13452.882

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 232.8
$ws.Range("J22").Value = 610
$ws.Range("L22").Value = 610
$ws.Range("N22").Value = -956

$ws.Range("H134").Value = 63060.89
$ws.Range("I134").Value = 101548.18
$ws.Range("J134").Value = 2580.8572
$ws.Range("K134").Value = 304644.54
$ws.Range("L134").Value = 7742.571599999999
$ws.Range("M134").Value = -302109.54
$ws.Range("N134").Value = -12812.5716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2183.5417
$ws.Range("J31").Value = 3951.625
$ws.Range("L31").Value = 3951.625
$ws.Range("N31").Value = -4541.625

$ws.Range("H34").Value = 2183.5417
$ws.Range("J34").Value = 3951.625
$ws.Range("L34").Value = 3951.625
$ws.Range("N34").Value = -4355.625

$ws.Range("H99").Value = 69762.5
$ws.Range("I99").Value = 32790.5
$ws.Range("J99").Value = 168354.5
$ws.Range("K99").Value = 32790.5
$ws.Range("L99").Value = 168354.5
$ws.Range("M99").Value = -31292.5
$ws.Range("N99").Value = -171350.5

$ws.Range("H122").Value = 3476239.8
$ws.Range("I122").Value = 4634414.5
$ws.Range("K122").Value = 13903243.5
$ws.Range("M122").Value = -13900793.5

$ws.Range("H126").Value = 69762.5
$ws.Range("I126").Value = 32790.5
$ws.Range("J126").Value = 168354.5
$ws.Range("K126").Value = 98371.5
$ws.Range("L126").Value = 505063.5
$ws.Range("M126").Value = -95901.5
$ws.Range("N126").Value = -510003.5

$ws.Range("H132").Value = 1534.841
$ws.Range("I132").Value = 1329.7805
$ws.Range("J132").Value = 4337.3335
$ws.Range("K132").Value = 3989.3415
$ws.Range("L132").Value = 13012.0005
$ws.Range("M132").Value = -1459.3415
$ws.Range("N132").Value = -18072.0005

$ws.Range("H134").Value = 2154.8372
$ws.Range("I134").Value = 1882.2572
$ws.Range("J134").Value = 3347.375
$ws.Range("K134").Value = 5646.7716
$ws.Range("L134").Value = 10042.125
$ws.Range("M134").Value = -3111.7716
$ws.Range("N134").Value = -15112.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2773.818
$ws.Range("I102").Value = 2666.6667
$ws.Range("J102").Value = 2902.4
$ws.Range("K102").Value = 2666.6667
$ws.Range("L102").Value = 2902.4
$ws.Range("M102").Value = -1044.6667
$ws.Range("N102").Value = -6146.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1476.6666
$ws.Range("I7").Value = 1623.7778
$ws.Range("J7").Value = 1035.3334
$ws.Range("K7").Value = 1623.7778
$ws.Range("L7").Value = 1035.3334
$ws.Range("M7").Value = -1511.7778
$ws.Range("N7").Value = -1259.3334

$ws.Range("H22").Value = 383.70587
$ws.Range("I22").Value = 425.14285
$ws.Range("K22").Value = 425.14285
$ws.Range("M22").Value = -130.14285

$ws.Range("H27").Value = 383.70587
$ws.Range("I27").Value = 425.14285
$ws.Range("K27").Value = 425.14285
$ws.Range("M27").Value = -318.14285

$ws.Range("H46").Value = 916
$ws.Range("I46").Value = 875.55554
$ws.Range("J46").Value = 1077.7778
$ws.Range("K46").Value = 875.55554
$ws.Range("L46").Value = 1077.7778
$ws.Range("M46").Value = -687.55554
$ws.Range("N46").Value = -1453.7778

$ws.Range("H100").Value = 2127.1428
$ws.Range("I100").Value = 2266.6667
$ws.Range("J100").Value = 2022.5
$ws.Range("K100").Value = 2266.6667
$ws.Range("L100").Value = 2022.5
$ws.Range("M100").Value = -1725.6667
$ws.Range("N100").Value = -3104.5

$ws.Range("H126").Value = 1476.6666
$ws.Range("I126").Value = 1623.7778
$ws.Range("J126").Value = 1035.3334
$ws.Range("K126").Value = 4871.3334
$ws.Range("L126").Value = 3106.0002
$ws.Range("M126").Value = -2401.3334
$ws.Range("N126").Value = -8046.0002

$ws.Range("H132").Value = 1821.575
$ws.Range("I132").Value = 1638.44
$ws.Range("J132").Value = 2126.8
$ws.Range("K132").Value = 4915.32
$ws.Range("L132").Value = 6380.400000000001
$ws.Range("M132").Value = -2385.32
$ws.Range("N132").Value = -11440.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2474.75
$ws.Range("I126").Value = 1864
$ws.Range("J126").Value = 6750
$ws.Range("K126").Value = 5592
$ws.Range("L126").Value = 20250
$ws.Range("M126").Value = -3122
$ws.Range("N126").Value = -25190

$ws.Range("H132").Value = 3592.5
$ws.Range("I132").Value = 4128.553
$ws.Range("K132").Value = 12385.659
$ws.Range("M132").Value = -9855.659

$ws.Range("H136").Value = 7960.8203
$ws.Range("I136").Value = 9991.741
$ws.Range("J136").Value = 3391.25
$ws.Range("K136").Value = 29975.223
$ws.Range("L136").Value = 10173.75
$ws.Range("M136").Value = -27425.223
$ws.Range("N136").Value = -15273.75
